# Auto-generated edit script applying the Sargatanas_Profits.xlsx diff
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 29.758621
$ws.Cells.Item(2, 9).Value = 26.962963
$ws.Cells.Item(2, 10).Value = 67.5
$ws.Cells.Item(2, 11).Value = 26.962963
$ws.Cells.Item(2, 12).Value = 67.5
$ws.Cells.Item(2, 13).Value = 86.037037
$ws.Cells.Item(2, 14).Value = -293.5
$ws.Cells.Item(6, 8).Value = 156.6
$ws.Cells.Item(6, 9).Value = 129.71428
$ws.Cells.Item(6, 11).Value = 389.14284
$ws.Cells.Item(6, 13).Value = -277.14284
$ws.Cells.Item(8, 8).Value = 271.43478
$ws.Cells.Item(8, 9).Value = 45.4
$ws.Cells.Item(8, 11).Value = 136.2
$ws.Cells.Item(8, 13).Value = 2.800000000000011
$ws.Cells.Item(39, 8).Value = 1478
$ws.Cells.Item(39, 9).Value = 542.9091
$ws.Cells.Item(39, 10).Value = 2947.4285
$ws.Cells.Item(39, 11).Value = 1628.7273
$ws.Cells.Item(39, 12).Value = 8842.2855
$ws.Cells.Item(39, 13).Value = -1332.7273
$ws.Cells.Item(39, 14).Value = -9434.2855
$ws.Cells.Item(111, 8).Value = 31251732
$ws.Cells.Item(111, 10).Value = 1964.5
$ws.Cells.Item(111, 12).Value = 5893.5
$ws.Cells.Item(111, 14).Value = -12027.5
$ws.Cells.Item(129, 8).Value = 1300.5264
$ws.Cells.Item(129, 9).Value = 781.1818
$ws.Cells.Item(129, 11).Value = 2343.5454
$ws.Cells.Item(129, 13).Value = 2656.4546
$ws.Cells.Item(135, 8).Value = 1538735.8
$ws.Cells.Item(135, 10).Value = 310
$ws.Cells.Item(135, 12).Value = 2790
$ws.Cells.Item(135, 14).Value = -7860
$ws.Cells.Item(137, 8).Value = 4399.2
$ws.Cells.Item(137, 9).Value = 4624
$ws.Cells.Item(137, 10).Value = 4142.2856
$ws.Cells.Item(137, 11).Value = 13872
$ws.Cells.Item(137, 12).Value = 12426.8568
$ws.Cells.Item(137, 13).Value = -11322
$ws.Cells.Item(137, 14).Value = -17526.8568
$ws.Cells.Item(138, 8).Value = 1474638
$ws.Cells.Item(138, 9).Value = 2034.8636
$ws.Cells.Item(138, 10).Value = 2178926.5
$ws.Cells.Item(138, 11).Value = 6104.5908
$ws.Cells.Item(138, 12).Value = 6536779.5
$ws.Cells.Item(138, 13).Value = -964.5908
$ws.Cells.Item(138, 14).Value = -6547059.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1902.3846
$ws.Cells.Item(45, 9).Value = 2018
$ws.Cells.Item(45, 11).Value = 2018
$ws.Cells.Item(45, 13).Value = -1641
$ws.Cells.Item(61, 8).Value = 25005194
$ws.Cells.Item(61, 9).Value = 2472.8
$ws.Cells.Item(61, 11).Value = 2472.8
$ws.Cells.Item(61, 13).Value = -2260.8
$ws.Cells.Item(132, 8).Value = 4605.646
$ws.Cells.Item(132, 9).Value = 2797.6123
$ws.Cells.Item(132, 11).Value = 8392.836899999998
$ws.Cells.Item(132, 13).Value = -5862.836899999998
$ws.Cells.Item(136, 8).Value = 25005194
$ws.Cells.Item(136, 9).Value = 2472.8
$ws.Cells.Item(136, 11).Value = 7418.400000000001
$ws.Cells.Item(136, 13).Value = -4868.400000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 54597.527
$ws.Cells.Item(86, 9).Value = 77814.305
$ws.Cells.Item(86, 10).Value = 4294.5
$ws.Cells.Item(86, 11).Value = 77814.305
$ws.Cells.Item(86, 12).Value = 4294.5
$ws.Cells.Item(86, 13).Value = -76691.305
$ws.Cells.Item(86, 14).Value = -6540.5
$ws.Cells.Item(89, 8).Value = 54597.527
$ws.Cells.Item(89, 9).Value = 77814.305
$ws.Cells.Item(89, 10).Value = 4294.5
$ws.Cells.Item(89, 11).Value = 389071.525
$ws.Cells.Item(89, 12).Value = 21472.5
$ws.Cells.Item(89, 13).Value = -383455.525
$ws.Cells.Item(89, 14).Value = -32704.5
$ws.Cells.Item(94, 8).Value = 948.54285
$ws.Cells.Item(94, 9).Value = 661.1613
$ws.Cells.Item(94, 11).Value = 661.1613
$ws.Cells.Item(94, 13).Value = -210.1613
$ws.Cells.Item(99, 8).Value = 5052584
$ws.Cells.Item(99, 9).Value = 1833.4166
$ws.Cells.Item(99, 10).Value = 15154086
$ws.Cells.Item(99, 11).Value = 1833.4166
$ws.Cells.Item(99, 12).Value = 15154086
$ws.Cells.Item(99, 13).Value = -335.4166
$ws.Cells.Item(99, 14).Value = -15157082
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6094.306
$ws.Cells.Item(31, 9).Value = 1677.9584
$ws.Cells.Item(31, 10).Value = 10334
$ws.Cells.Item(31, 11).Value = 1677.9584
$ws.Cells.Item(31, 12).Value = 10334
$ws.Cells.Item(31, 13).Value = -1382.9584
$ws.Cells.Item(31, 14).Value = -10924
$ws.Cells.Item(34, 8).Value = 6094.306
$ws.Cells.Item(34, 9).Value = 1677.9584
$ws.Cells.Item(34, 10).Value = 10334
$ws.Cells.Item(34, 11).Value = 1677.9584
$ws.Cells.Item(34, 12).Value = 10334
$ws.Cells.Item(34, 13).Value = -1475.9584
$ws.Cells.Item(34, 14).Value = -10738
$ws.Cells.Item(58, 8).Value = 6290.2856
$ws.Cells.Item(58, 9).Value = 1841.6666
$ws.Cells.Item(58, 10).Value = 8611.305
$ws.Cells.Item(58, 11).Value = 1841.6666
$ws.Cells.Item(58, 12).Value = 8611.305
$ws.Cells.Item(58, 13).Value = -1638.6666
$ws.Cells.Item(58, 14).Value = -9017.305
$ws.Cells.Item(69, 8).Value = 27000
$ws.Cells.Item(69, 9).Value = 27000
$ws.Cells.Item(69, 11).Value = 27000
$ws.Cells.Item(69, 13).Value = -26251
$ws.Cells.Item(72, 8).Value = 27000
$ws.Cells.Item(72, 9).Value = 27000
$ws.Cells.Item(72, 11).Value = 81000
$ws.Cells.Item(72, 13).Value = -77256
$ws.Cells.Item(86, 8).Value = 4176772.5
$ws.Cells.Item(86, 10).Value = 7715.6665
$ws.Cells.Item(86, 12).Value = 7715.6665
$ws.Cells.Item(86, 14).Value = -9961.6665
$ws.Cells.Item(89, 8).Value = 4176772.5
$ws.Cells.Item(89, 10).Value = 7715.6665
$ws.Cells.Item(89, 12).Value = 38578.3325
$ws.Cells.Item(89, 14).Value = -49810.3325
$ws.Cells.Item(136, 8).Value = 6290.2856
$ws.Cells.Item(136, 9).Value = 1841.6666
$ws.Cells.Item(136, 10).Value = 8611.305
$ws.Cells.Item(136, 11).Value = 5524.9998
$ws.Cells.Item(136, 12).Value = 25833.915
$ws.Cells.Item(136, 13).Value = -2974.9998
$ws.Cells.Item(136, 14).Value = -30933.915
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 78339.56
$ws.Cells.Item(2, 9).Value = 13684.954
$ws.Cells.Item(2, 10).Value = 220579.7
$ws.Cells.Item(2, 11).Value = 82109.724
$ws.Cells.Item(2, 12).Value = 1323478.2
$ws.Cells.Item(2, 13).Value = -81996.724
$ws.Cells.Item(2, 14).Value = -1323704.2
$ws.Cells.Item(44, 8).Value = 1060
$ws.Cells.Item(44, 10).Value = 1137.5
$ws.Cells.Item(44, 12).Value = 3412.5
$ws.Cells.Item(44, 14).Value = -4208.5
$ws.Cells.Item(56, 8).Value = 8000
$ws.Cells.Item(56, 9).Value = 8000
$ws.Cells.Item(56, 11).Value = 8000
$ws.Cells.Item(56, 13).Value = -7470
$ws.Cells.Item(62, 8).Value = 3142.9285
$ws.Cells.Item(62, 10).Value = 3153.923
$ws.Cells.Item(62, 12).Value = 9461.769
$ws.Cells.Item(62, 14).Value = -10833.769
$ws.Cells.Item(65, 8).Value = 3142.9285
$ws.Cells.Item(65, 10).Value = 3153.923
$ws.Cells.Item(65, 12).Value = 28385.307
$ws.Cells.Item(65, 14).Value = -35249.307
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).ClearContents()
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(127, 8).Value = 5044
$ws.Cells.Item(127, 10).Value = 5044
$ws.Cells.Item(127, 12).Value = 15132
$ws.Cells.Item(127, 14).Value = -25052
$ws.Cells.Item(141, 8).Value = 4146.6
$ws.Cells.Item(141, 9).Value = 4146.6
$ws.Cells.Item(141, 11).Value = 12439.8
$ws.Cells.Item(141, 13).Value = -7259.800000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 927.2273
$ws.Cells.Item(97, 9).Value = 1051.75
$ws.Cells.Item(97, 11).Value = 1051.75
$ws.Cells.Item(97, 13).Value = -555.75
$ws.Cells.Item(107, 8).Value = 458.06668
$ws.Cells.Item(107, 9).Value = 624.625
$ws.Cells.Item(107, 11).Value = 624.625
$ws.Cells.Item(107, 13).Value = 1295.375
$ws.Cells.Item(113, 8).Value = 208491.06
$ws.Cells.Item(113, 9).Value = 436444.88
$ws.Cells.Item(113, 11).Value = 436444.88
$ws.Cells.Item(113, 13).Value = -434274.88
$ws.Cells.Item(132, 8).Value = 2661.4614
$ws.Cells.Item(132, 9).Value = 1410.591
$ws.Cells.Item(132, 11).Value = 4231.772999999999
$ws.Cells.Item(132, 13).Value = -1701.772999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4711.381
$ws.Cells.Item(7, 9).Value = 3576.3333
$ws.Cells.Item(7, 10).Value = 5562.6665
$ws.Cells.Item(7, 11).Value = 3576.3333
$ws.Cells.Item(7, 12).Value = 5562.6665
$ws.Cells.Item(7, 13).Value = -3464.3333
$ws.Cells.Item(7, 14).Value = -5786.6665
$ws.Cells.Item(22, 8).Value = 3779
$ws.Cells.Item(22, 9).Value = 1985.7142
$ws.Cells.Item(22, 11).Value = 1985.7142
$ws.Cells.Item(22, 13).Value = -1690.7142
$ws.Cells.Item(27, 8).Value = 3779
$ws.Cells.Item(27, 9).Value = 1985.7142
$ws.Cells.Item(27, 11).Value = 1985.7142
$ws.Cells.Item(27, 13).Value = -1878.7142
$ws.Cells.Item(40, 8).Value = 5070.7856
$ws.Cells.Item(40, 9).Value = 2141.7144
$ws.Cells.Item(40, 11).Value = 2141.7144
$ws.Cells.Item(40, 13).Value = -2005.7144
$ws.Cells.Item(126, 8).Value = 4711.381
$ws.Cells.Item(126, 9).Value = 3576.3333
$ws.Cells.Item(126, 10).Value = 5562.6665
$ws.Cells.Item(126, 11).Value = 10728.9999
$ws.Cells.Item(126, 12).Value = 16687.9995
$ws.Cells.Item(126, 13).Value = -8258.999899999999
$ws.Cells.Item(126, 14).Value = -21627.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).ClearContents()
$ws.Cells.Item(96, 14).ClearContents()
